$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsMeta = $wb.Worksheets.Item(1)
$wsConf = $wb.Worksheets.Item(2)   # currently "Include from Match Confidence"
$wsGrade = $wb.Worksheets.Item(3)  # currently "Include from match-grade"

# --- Metadata sheet updates ---
# Version 5.0.0 -> 6.0.0
$wsMeta.Range("B3").Value = "6.0.0"

# Date updated
$wsMeta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Row 9 was Publisher | (blank) -> Publisher | Alvearie Team
$wsMeta.Range("B9").Value = "Alvearie Team"

# Row 10 was Contact | No display for ContactDetail -> Jurisdiction | United States of America
$wsMeta.Range("A10").Value = "Jurisdiction"
$wsMeta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate Contact | No display for ContactDetail row - remove it by
# deleting the whole row, which shifts rows 12-15 up to 11-14
$wsMeta.Rows.Item(11).Delete()

# --- Rename the two "Include from ..." sheets and swap their System URI content ---
$wsConf.Name = "Include from MatchGrade"
$wsGrade.Name = "Include from Match Confidence"

# The (now renamed) MatchGrade sheet's System URI value
$wsConf.Range("B4").Value = "http://terminology.hl7.org/CodeSystem/match-grade"

# The (now renamed) Match Confidence sheet's System URI value
$wsGrade.Range("B4").Value = "http://ibm.com/fhir/cdm/CodeSystem/match-confidence-level"
